# ------------------------------------------------------------------
# Applies the NH3 optimisation re-run numbers (commit: "added
# project10 and other fixes") to nh3_opt.xlsx's three sheets.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ====================================================================
# Sheet "Energies"
# ====================================================================
$wsE = $wb.Worksheets.Item("Energies")

# Total row
$wsE.Range("B2").Value = 23.376
$wsE.Range("C2").Value = 6.424
$wsE.Range("D2").Value = 48.178

# Rotational row
$wsE.Range("D5").Value = 13.631

# Vibrational row
$wsE.Range("B6").Value = 21.598
$wsE.Range("C6").Value = 0.463
$wsE.Range("D6").Value = 0.106

# ====================================================================
# Sheet "Partition Functions"
# ====================================================================
$wsP = $wb.Worksheets.Item("Partition Functions")

# Rotational row
$wsP.Range("B4").Value = 212.608
$wsP.Range("C4").Value = 2.327579
$wsP.Range("D4").Value = 5.35945

# Total Bot row
$wsP.Range("B5").Value = [double]"9.12601e-08"
$wsP.Range("C5").Value = -7.039719
$wsP.Range("D5").Value = -16.209552

# Total V=0 row
$wsP.Range("B6").Value = [double]"592331000"
$wsP.Range("C6").Value = [double]"8.772563999999999"
$wsP.Range("D6").Value = 20.199576

# Vib (Bot) row
$wsP.Range("B7").Value = [double]"1.55438e-16"
$wsP.Range("C7").Value = -15.808444
$wsP.Range("D7").Value = -36.400287

# Vib (V=0) row
$wsP.Range("B8").Value = 1.00888
$wsP.Range("C8").Value = 0.00384
$wsP.Range("D8").Value = [double]"0.008841999999999999"

# ====================================================================
# Sheet "Other"
# ====================================================================
$wsO = $wb.Worksheets.Item("Other")

# The duplicated last SCF energy row (old row 6) is dropped, which
# shifts every following block up by one row.
$wsO.Rows(6).Delete()

# A brand-new vibrational temperature entry appears (the former
# imaginary-frequency mode now has a positive vibrational
# temperature), which pushes the remaining five values back down by
# one row.
$wsO.Rows(35).Insert()

# --- Energies SCF [ Hartree ] -------------------------------------
$wsO.Range("A2").Value = -56.5662326571
$wsO.Range("A3").Value = -56.5669741282
$wsO.Range("A4").Value = -56.5669844677
$wsO.Range("A5").Value = -56.5669844677

# --- Temperature [ K ] ---------------------------------------------
$wsO.Range("A8").Value = 298.15

# --- Pressure [ atm ] ------------------------------------------------
$wsO.Range("A11").Value = 1

# --- Principal Moments of Inertia [ amu Å^2 ] -----------------------
$wsO.Range("A14").Value = 6.01866
$wsO.Range("B14").Value = 6.01934
$wsO.Range("C14").Value = 9.73682

# --- Molecular Mass [ amu ] -----------------------------------------
$wsO.Range("A17").Value = 17.02655

# --- Rotational Symmetry Number -------------------------------------
$wsO.Range("A20").Value = 1

# --- Rotational Temperatures [ K ] ----------------------------------
$wsO.Range("A23").Value = 14.39089
$wsO.Range("B23").Value = 14.38924
$wsO.Range("C23").Value = [double]"8.895490000000001"

# --- Rotational Constants [ GHZ ] -----------------------------------
$wsO.Range("A26").Value = 299.85786
$wsO.Range("B26").Value = 299.82356
$wsO.Range("C26").Value = 185.35216

# --- Zero-point Vibrational Energy [ J/mol ] ------------------------
$wsO.Range("A29").Value = 90256.8

# --- Zero-point Vibrational Energy [ Kcal/mol ] ---------------------
$wsO.Range("A32").Value = 21.57189

# --- Vibrational Temperatures [ K ] ---------------------------------
$wsO.Range("A35").Value = 1432.76
$wsO.Range("A36").Value = 2407.14
$wsO.Range("A37").Value = 2407.46
$wsO.Range("A38").Value = 5016.86
$wsO.Range("A39").Value = 5222.94
$wsO.Range("A40").Value = 5223.61

# --- Zero-point Correction [ Hartree/Particle ] ---------------------
$wsO.Range("A43").Value = 0.034377

# --- Thermal Correction to Energy [ Hartree/Particle ] --------------
$wsO.Range("A46").Value = 0.037252

# --- Thermal Correction to Enthalpy [ Hartree/Particle ] ------------
$wsO.Range("A49").Value = 0.038196

# --- Thermal Correction to Gibbs Free Energy [ Hartree/Particle ] ---
$wsO.Range("A52").Value = 0.015305

# --- Harmonic Frequencies [ cm^-1 ] ----------------------------------
$wsO.Range("A55").Value = 995.8161
$wsO.Range("A56").Value = 1673.0462
$wsO.Range("A57").Value = 1673.2702
$wsO.Range("A58").Value = 3486.8987
$wsO.Range("A59").Value = 3630.1303
$wsO.Range("A60").Value = 3630.5925

# --- Reduced Masses [ amu ] ------------------------------------------
$wsO.Range("A63").Value = 1.1841
$wsO.Range("A64").Value = 1.0685
$wsO.Range("A65").Value = 1.0685
$wsO.Range("A66").Value = 1.0241
$wsO.Range("A67").Value = 1.0913
$wsO.Range("A68").Value = 1.0913

# --- Force Constants [ mDyne/A ] -------------------------------------
$wsO.Range("A71").Value = 0.6918
$wsO.Range("A72").Value = 1.7621
$wsO.Range("A73").Value = 1.7626
$wsO.Range("A74").Value = 7.3363
$wsO.Range("A75").Value = 8.4732
$wsO.Range("A76").Value = 8.4756

# --- IR Intensities [ km/mol ] ---------------------------------------
$wsO.Range("A79").Value = 241.5389
$wsO.Range("A80").Value = 28.7767
$wsO.Range("A81").Value = 28.7711
$wsO.Range("A82").Value = 2.6039
$wsO.Range("A83").Value = 3.9852
$wsO.Range("A84").Value = 4.0005

Write-Host "done"
